# Apply trade #26 close-out update to the live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - update headline stats
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.15
$wsSummary.Range("B4").Value = -0.86
$wsSummary.Range("B5").Value = -0.66
$wsSummary.Range("B6").Value = 26
$wsSummary.Range("B8").Value = 16
$wsSummary.Range("B9").Value = 30.77

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - update MarketMaking row (row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.15000000000001
$wsStatus.Range("D4").Value = 26
$wsStatus.Range("E4").Value = -0.86
$wsStatus.Range("F4").Value = -0.85
$wsStatus.Range("G4").Value = 30.77

# ---------------------------------------------------------------------
# 3) Append new trade #26 row (row 27) to "All Trades" and "MarketMaking"
#    sheets. The "Date" column (B) looks like a date to Excel's
#    auto-detection, so it is entered with a leading apostrophe to force
#    it to stay literal text (matching the source data, which stores it
#    as a plain string rather than a date serial). The other text
#    columns (time/strategy/side/status/reason) already round-trip as
#    plain strings without any extra hinting.
# ---------------------------------------------------------------------
function Add-TradeRow($ws) {
    $row = 27

    $ws.Cells.Item($row, 1).Value = 26
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "13:18:53"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.16
    $ws.Cells.Item($row, 7).Value = 0.11
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -31.25
    $ws.Cells.Item($row, 10).Value = -0.05
    $ws.Cells.Item($row, 11).Value = 99.15000000000001
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
